$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.174752354621887
$ws.Range("B1").Value = 1.426955223083496
$ws.Range("C1").Value = 1.774097442626953
$ws.Range("D1").Value = 1.618454813957214
$ws.Range("E1").Value = 1.563287615776062
